$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("final")
$lo = $ws.ListObjects.Item(4)
Write-Host ("Name: " + $lo.Name + " Range: " + $lo.Range.Address())
for ($i=1; $i -le $lo.ListColumns.Count; $i++) {
    $col = $lo.ListColumns.Item($i)
    Write-Host ("  " + $i.ToString() + ": " + $col.Name + " body=" + $col.DataBodyRange.Address())
}
